# Update "想去人数" (interest count) figures for several events that are
# shared across the "展览" (Exhibition) sheet and the "全部类型" (All types)
# sheet, reflecting newly-generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 5345
$ws1.Range("F8").Value  = 880
$ws1.Range("F9").Value  = 130
$ws1.Range("F10").Value = 2373
$ws1.Range("F12").Value = 53
$ws1.Range("F13").Value = 2224

# Sheet "全部类型": same events, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 5345
$ws4.Range("F10").Value = 880
$ws4.Range("F11").Value = 130
$ws4.Range("F12").Value = 2373
$ws4.Range("F15").Value = 53
$ws4.Range("F16").Value = 2224

$wb.Save()
